$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LISTE")
$lo = $ws.ListObjects.Item("Tableau13")

# Add four new rows to the "Tableau13" table (expands the table range
# automatically from A5:M86 to A5:M90).
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# Row 87 - Diode TVS (Comchip)
$ws.Range("A87").Value = "E00077"
$ws.Range("B87").Value = "Composant"
$ws.Range("C87").Value = "Diodes TVS"
$ws.Range("D87").Value = "CPDQ5V0-HF"
$ws.Range("E87").Value = "Diode TVS "
$ws.Range("F87").Value = "Comchip"
$ws.Range("G87").Value = "E00077"
$ws.Range("H87").Value = "V1.00"

# Row 88 - Polyswitch (Bel Fuse Inc)
$ws.Range("A88").Value = "E00081"
$ws.Range("B88").Value = "Composant"
$ws.Range("C88").Value = "Fuse"
$ws.Range("D88").Value = "0ZCJ0100FF2E"
$ws.Range("E88").Value = "Polyswitch"
$ws.Range("F88").Value = "Bel Fuse Inc"
$ws.Range("G88").Value = "E00081"
$ws.Range("H88").Value = "V1.00"
$ws.Range("L88").Value = "0ZCJ0100FF2E"

# Row 89 - Transistor (Mosfet P 20V 4.4A, Toshiba)
$ws.Range("A89").Value = "E00082"
$ws.Range("B89").Value = "Composant"
$ws.Range("C89").Value = "Transistor "
$ws.Range("D89").Value = "SSM3J130TU"
$ws.Range("E89").Value = "Mosfet P 20V 4.4A"
$ws.Range("F89").Value = "Toshiba"
$ws.Range("G89").Value = "E00082"
$ws.Range("H89").Value = "V1.00"
$ws.Range("L89").Value = "SSM3J130TU"

# Row 90 - Transistor (Bipolaire NPN)
$ws.Range("A90").Value = "E00083"
$ws.Range("B90").Value = "Composant"
$ws.Range("C90").Value = "Transistor "
$ws.Range("D90").Value = "MMBT2222ATT1G"
$ws.Range("E90").Value = "Bipolaire NPN "
$ws.Range("G90").Value = "E00083"
$ws.Range("H90").Value = "V1.00"
$ws.Range("L90").Value = "MMBT2222ATT1G"

# Hyperlink on the new REPERTOIRE cell, consistent with the other rows
# in the column (each links to its own external "E000xx" reference).
$ws.Hyperlinks.Add($ws.Range("G87"), "E00077") | Out-Null

# Keep the view roughly where the new rows were edited, as in the source file.
$ws.Range("F90").Select() | Out-Null
